$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# This document has two identical short list headings with the text
# "Часть 1" (bold, Times New Roman). The first one belongs to the list item
# introducing Part 1 and must stay untouched; the second one (which follows
# the screenshot/picture illustrating Part 1's steps) is mislabeled and
# should read "Часть 2" instead, since it introduces the second part of the
# assignment. Word also relocates its auto-managed "_GoBack" bookmark
# (marking the last edit position) to the point right after this edit.
# ---------------------------------------------------------------------------

# Drop the existing "_GoBack" bookmark - it will be re-created at the new
# edit location below (Word always keeps at most one "_GoBack" bookmark).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Locate the target paragraph: the short "Часть 1" list heading that comes
# after the first inline picture in the document (i.e. the second, later
# occurrence of this heading - the one that needs to become "Часть 2").
$picEnd = 0
if ($d.InlineShapes.Count -gt 0) {
    $picEnd = $d.InlineShapes(1).Range.End
}

$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -eq "Часть 1`r") {
        if ($p.Range.Start -gt $picEnd) {
            $target = $p
            break
        }
    }
}

$pEnd = $target.Range.End

# Replace the trailing "1" with "2", keeping it as its own run (matching
# how Word splits a run when only part of its text is retyped) rather than
# merging back into the preceding "Часть " run.
$oldDigit = $d.Range($pEnd - 2, $pEnd - 1)
$oldDigit.Delete()

$insPt = $d.Range($pEnd - 2, $pEnd - 2)
$insPt.InsertAfter("2")

$twoRange = $d.Range($pEnd - 2, $pEnd - 1)
$twoRange.Font.Bold = 0
$twoRange.Font.Bold = 1

# Re-establish the "_GoBack" bookmark right after the edit, at the end of
# the paragraph (a zero-width "point" bookmark, as Word leaves it).
$pEnd2 = $target.Range.End
$placeholder = $d.Range($pEnd2 - 1, $pEnd2 - 1)
$placeholder.InsertAfter("X")
$wrapRange = $d.Range($pEnd2 - 1, $pEnd2)
$d.Bookmarks.Add("_GoBack", $wrapRange)
$d.Range($pEnd2 - 1, $pEnd2).Delete()
